$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the package id (FInterID) from 6500 to 6507 for the existing BOM rows (2-7)
$ws.Range("A2:A7").Value = 6507

# The last line item (row 8, FEntryID 7) is no longer part of the finished BOM - remove it entirely
$ws.Rows("8:8").Delete()
